$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 551.5
$ws.Range("I19").Value = 270
$ws.Range("J19").Value = 591.7143
$ws.Range("K19").Value = 270
$ws.Range("L19").Value = 591.7143
$ws.Range("M19").Value = -95
$ws.Range("N19").Value = -941.7143
$ws.Range("H33").Value = 6697.8076
$ws.Range("I33").Value = 6697.8076
$ws.Range("K33").Value = 6697.8076
$ws.Range("M33").Value = -6468.8076
$ws.Range("H41").Value = 1470.5
$ws.Range("I41").Value = 670
$ws.Range("J41").Value = 1901.5385
$ws.Range("K41").Value = 670
$ws.Range("L41").Value = 1901.5385
$ws.Range("M41").Value = -230
$ws.Range("N41").Value = -2781.5385
$ws.Range("H82").Value = 1823.3334
$ws.Range("I82").Value = 1823.3334
$ws.Range("K82").Value = 5470.0002
$ws.Range("M82").Value = -5064.0002
$ws.Range("H85").Value = 1823.3334
$ws.Range("I85").Value = 1823.3334
$ws.Range("K85").Value = 5470.0002
$ws.Range("M85").Value = -4066.0002
$ws.Range("H100").Value = 2599.25
$ws.Range("I100").Value = 1199.5
$ws.Range("K100").Value = 1199.5
$ws.Range("M100").Value = -658.5
$ws.Range("H112").Value = 334939.94
$ws.Range("I112").Value = 1114
$ws.Range("J112").Value = 436539.12
$ws.Range("K112").Value = 3342
$ws.Range("L112").Value = 1309617.36
$ws.Range("M112").Value = -2234
$ws.Range("N112").Value = -1311833.36
$ws.Range("H113").Value = 9610.825999999999
$ws.Range("I113").Value = 10540.385
$ws.Range("J113").Value = 8402.4
$ws.Range("K113").Value = 10540.385
$ws.Range("L113").Value = 8402.4
$ws.Range("M113").Value = -7286.385
$ws.Range("N113").Value = -14910.4
$ws.Range("H137").Value = 434991.25
$ws.Range("I137").Value = 1565.1765
$ws.Range("K137").Value = 4695.529500000001
$ws.Range("M137").Value = -2145.529500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 1430.5
$ws.Range("J29").Value = 1430.5
$ws.Range("L29").Value = 1430.5
$ws.Range("N29").Value = -2046.5
$ws.Range("H32").Value = 24933.666
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 24933.666
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 24933.666
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -25507.666
$ws.Range("H74").Value = 3990.5
$ws.Range("I74").Value = 1984
$ws.Range("J74").Value = 8003.5
$ws.Range("K74").Value = 1984
$ws.Range("L74").Value = 8003.5
$ws.Range("M74").Value = -1110
$ws.Range("N74").Value = -9751.5
$ws.Range("H77").Value = 3990.5
$ws.Range("I77").Value = 1984
$ws.Range("J77").Value = 8003.5
$ws.Range("K77").Value = 9920
$ws.Range("L77").Value = 40017.5
$ws.Range("M77").Value = -5552
$ws.Range("N77").Value = -48753.5
$ws.Range("H102").Value = 43971.43
$ws.Range("I102").Value = 21555.8
$ws.Range("K102").Value = 21555.8
$ws.Range("M102").Value = -19933.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 102292.37
$ws.Range("J20").Value = 2822.4285
$ws.Range("L20").Value = 2822.4285
$ws.Range("N20").Value = -3316.4285
$ws.Range("H94").Value = 830.34375
$ws.Range("I94").Value = 746.75
$ws.Range("K94").Value = 746.75
$ws.Range("M94").Value = -295.75
$ws.Range("H105").Value = 2399.75
$ws.Range("I105").Value = 2116.3333
$ws.Range("K105").Value = 2116.3333
$ws.Range("M105").Value = -369.3332999999998
$ws.Range("H124").Value = 69999.5
$ws.Range("J124").Value = 69999.5
$ws.Range("L124").Value = 69999.5
$ws.Range("N124").Value = -79819.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3743.718
$ws.Range("I31").Value = 2445.2083
$ws.Range("J31").Value = 5821.3335
$ws.Range("K31").Value = 2445.2083
$ws.Range("L31").Value = 5821.3335
$ws.Range("M31").Value = -2150.2083
$ws.Range("N31").Value = -6411.3335
$ws.Range("H34").Value = 3743.718
$ws.Range("I34").Value = 2445.2083
$ws.Range("J34").Value = 5821.3335
$ws.Range("K34").Value = 2445.2083
$ws.Range("L34").Value = 5821.3335
$ws.Range("M34").Value = -2243.2083
$ws.Range("N34").Value = -6225.3335
$ws.Range("H141").Value = 233726.1
$ws.Range("J141").Value = 254598.7
$ws.Range("L141").Value = 254598.7
$ws.Range("N141").Value = -264958.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 950
$ws.Range("I6").Value = 1000
$ws.Range("K6").Value = 3000
$ws.Range("M6").Value = -2887
$ws.Range("H86").Value = 2375.8386
$ws.Range("I86").Value = 355.1
$ws.Range("J86").Value = 3338.0952
$ws.Range("K86").Value = 1065.3
$ws.Range("L86").Value = 10014.2856
$ws.Range("M86").Value = 120.6999999999998
$ws.Range("N86").Value = -12386.2856
$ws.Range("H89").Value = 2375.8386
$ws.Range("I89").Value = 355.1
$ws.Range("J89").Value = 3338.0952
$ws.Range("K89").Value = 3195.9
$ws.Range("L89").Value = 30042.8568
$ws.Range("M89").Value = 2732.1
$ws.Range("N89").Value = -41898.8568
$ws.Range("H113").Value = 825.65515
$ws.Range("I113").Value = 552.4545000000001
$ws.Range("J113").Value = 992.6111
$ws.Range("K113").Value = 1657.3635
$ws.Range("L113").Value = 2977.8333
$ws.Range("M113").Value = 512.6364999999998
$ws.Range("N113").Value = -7317.8333
$ws.Range("H114").Value = 7565.0625
$ws.Range("I114").Value = 257
$ws.Range("J114").Value = 10001.083
$ws.Range("K114").Value = 771
$ws.Range("L114").Value = 30003.249
$ws.Range("M114").Value = 2483
$ws.Range("N114").Value = -36511.249
$ws.Range("H122").Value = 1123077.6
$ws.Range("J122").Value = 1263337.4
$ws.Range("L122").Value = 11370036.6
$ws.Range("N122").Value = -11374936.6
$ws.Range("H132").Value = 3506.7297
$ws.Range("I132").Value = 2330.75
$ws.Range("J132").Value = 3831.138
$ws.Range("K132").Value = 20976.75
$ws.Range("L132").Value = 34480.242
$ws.Range("M132").Value = -18446.75
$ws.Range("N132").Value = -39540.242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 478324.66
$ws.Range("I70").Value = 481099.78
$ws.Range("J70").Value = 469999.34
$ws.Range("K70").Value = 481099.78
$ws.Range("L70").Value = 469999.34
$ws.Range("M70").Value = -480829.78
$ws.Range("N70").Value = -470539.34
$ws.Range("H73").Value = 478324.66
$ws.Range("I73").Value = 481099.78
$ws.Range("J73").Value = 469999.34
$ws.Range("K73").Value = 481099.78
$ws.Range("L73").Value = 469999.34
$ws.Range("M73").Value = -480163.78
$ws.Range("N73").Value = -471871.34
$ws.Range("H113").Value = 2651511.5
$ws.Range("I113").Value = 222932
$ws.Range("K113").Value = 222932
$ws.Range("M113").Value = -220762
$ws.Range("H132").Value = 2182.9092
$ws.Range("I132").Value = 1770.25
$ws.Range("J132").Value = 3283.3333
$ws.Range("K132").Value = 5310.75
$ws.Range("L132").Value = 9849.999899999999
$ws.Range("M132").Value = -2780.75
$ws.Range("N132").Value = -14909.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 768.4286
$ws.Range("I16").Value = 768.4286
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 768.4286
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -598.4286
$ws.Range("N16").ClearContents()
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H55").Value = 2943227.8
$ws.Range("I55").Value = 1114.5834
$ws.Range("K55").Value = 1114.5834
$ws.Range("M55").Value = -941.5834
$ws.Range("H68").Value = 212595.2
$ws.Range("I68").Value = 264744
$ws.Range("K68").Value = 264744
$ws.Range("M68").Value = -263995
$ws.Range("H71").Value = 212595.2
$ws.Range("I71").Value = 264744
$ws.Range("K71").Value = 1323720
$ws.Range("M71").Value = -1319976
$ws.Range("H93").Value = 1834.3103
$ws.Range("I93").Value = 1536.7826
$ws.Range("J93").Value = 2974.8333
$ws.Range("K93").Value = 1536.7826
$ws.Range("L93").Value = 2974.8333
$ws.Range("M93").Value = -288.7826
$ws.Range("N93").Value = -5470.8333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30000
$ws.Range("I40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("K40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("M40").Value = -29851
$ws.Range("N40").Value = -30298
$ws.Range("H122").Value = 2900.8462
$ws.Range("I122").Value = 1540.5769
$ws.Range("K122").Value = 4621.7307
$ws.Range("M122").Value = -2171.7307
$ws.Range("H126").Value = 38271
$ws.Range("I126").Value = 56407.2
$ws.Range("J126").Value = 1998.6
$ws.Range("K126").Value = 169221.6
$ws.Range("L126").Value = 5995.799999999999
$ws.Range("M126").Value = -166751.6
$ws.Range("N126").Value = -10935.8
$ws.Range("H140").Value = 99473.75
$ws.Range("J140").Value = 99473.75
$ws.Range("L140").Value = 99473.75
$ws.Range("N140").Value = -109833.75
